$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.405.60"
$ws.Range("E2").Value = "  +1.11%  "
$ws.Range("D3").Value = "3.434.91"
$ws.Range("E3").Value = "  +1.92%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "575.25"
$ws.Range("E5").Value = "  +1.10%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.33"
$ws.Range("E6").Value = "  +7.15%  "
$ws.Range("D7").Value = "3.435.67"
$ws.Range("E7").Value = "  +1.99%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.479"
$ws.Range("E9").Value = "  +2.31%  "
$ws.Range("E10").Value = "  +1.11%  "
$ws.Range("E11").Value = "  +3.67%  "
$ws.Range("E12").Value = "  +2.20%  "
$ws.Range("D13").Value = "4.022.73"
$ws.Range("E13").Value = "  +1.97%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.98"
$ws.Range("E14").Value = "  +7.75%  "
$ws.Range("E15").Value = "  -0.59%  "
$ws.Range("E16").Value = "  +2.17%  "
$ws.Range("D17").Value = "3.436.31"
$ws.Range("E17").Value = "  +1.90%  "
$ws.Range("D18").Value = "61.512.19"
$ws.Range("E18").Value = "  +1.16%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.28"
$ws.Range("E19").Value = "  +8.23%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.21"
$ws.Range("E20").Value = "  +3.59%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.41"
$ws.Range("E21").Value = "  +2.30%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "395.66"
$ws.Range("E22").Value = "  +6.69%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.566"
$ws.Range("E23").Value = "  +3.60%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "73.41"
$ws.Range("E24").Value = "  +3.92%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.998"
$ws.Range("E25").Value = "  -0.35%  "
$ws.Range("E26").Value = "  +0.38%  "
$ws.Range("E27").Value = "  +0.33%  "
$ws.Range("D28").Value = "3.571.95"
$ws.Range("E28").Value = "  +1.81%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.179"
$ws.Range("E29").Value = "  +2.85%  "
$ws.Range("E30").Value = "  +3.59%  "
$ws.Range("E31").Value = "  +0.06%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.20"
$ws.Range("E32").Value = "  +2.29%  "
$ws.Range("E33").Value = "  -6.38%  "
$ws.Range("E34").Value = "  +2.59%  "
$ws.Range("E35").Value = "  +0.00%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "24.04"
$ws.Range("E36").Value = "  +3.50%  "
$ws.Range("B37").Value = "RenzoRestakedETH"
$ws.Range("C37").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D37").Value = "3.462.01"
$ws.Range("E37").Value = "  +2.19%  "
$ws.Range("B38").Value = "Aptos"
$ws.Range("C38").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "7.01"
$ws.Range("E38").Value = "  +3.93%  "
$ws.Range("E39").Value = "  +1.99%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.12"
$ws.Range("E40").Value = "  +0.60%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "167.57"
$ws.Range("E41").Value = "  +1.94%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0785"
$ws.Range("E42").Value = "  +3.47%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "26.80"
$ws.Range("E43").Value = "  +7.34%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.799"
$ws.Range("E44").Value = "  +3.86%  "
$ws.Range("B45").Value = "Stacks"
$ws.Range("C45").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.74"
$ws.Range("E45").Value = "  -0.06%  "
$ws.Range("B46").Value = "FirstDigitalUSD"
$ws.Range("C46").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.00"
$ws.Range("E46").Value = "  +0.00%  "
$ws.Range("E47").Value = "  +3.92%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "42.14"
$ws.Range("E48").Value = "  +0.66%  "
$ws.Range("B49").Value = "Maker"
$ws.Range("C49").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D49").Value = "2.588.69"
$ws.Range("E49").Value = "  +3.19%  "
$ws.Range("B50").Value = "ONDO"
$ws.Range("C50").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.16"
$ws.Range("E50").Value = "  +0.61%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.93"
$ws.Range("E51").Value = "  +2.83%  "
